# Apply updated dSF (column F) values per repull of data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 4
    4  = -1
    5  = -2
    6  = -4
    8  = -3
    9  = 1
    11 = 2
    17 = -3
    18 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
